# Commit: fix(case_csv): Replace dem with filled_dem
#
# The "Envs" worksheet contains many cells whose text is exactly "dem" or
# "dem " (trailing space). These need to become "filled_dem" and
# "filled_dem " respectively. Because these repeated text values live in
# the shared string table, using a Find/Replace across the sheet updates
# every cell that shares that string in one shot (matching how Excel's
# own Ctrl+H dialog behaves), rather than touching each cell individually.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Envs")
$ws.Activate()

# Excel's built-in constants for Find/Replace.
$xlWhole = 1          # xlWhole
$xlByRows = 1          # xlByRows
$xlNext = 1          # xlNext
$xlFormulas = -4123       # xlFormulas

# Replace the plain "dem" value first (matches whole-cell text only, so
# it cannot touch the "dem " cell), then the "dem " value with trailing
# space. Doing them in this order keeps the shared-string table entries
# in their original slots (index 241 stays "filled_dem", 242 stays
# "filled_dem ").
$ws.Cells.Replace("dem", "filled_dem", $xlWhole, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("dem ", "filled_dem ", $xlWhole, $xlByRows, $false, $false, $true, $true)

# Mirror the author's resulting cursor position: after the replace the
# first changed cell (C3) ends up selected/top-left in the saved file.
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 3
